$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains its original text formatting so that
# numeric-looking values (e.g. "4.31", "1.00") are not coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.810.16'
$ws.Range("E2").Value = '  -0.96%  '
$ws.Range("D3").Value = '1.614.08'
$ws.Range("E3").Value = '  -0.90%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '211.87'
$ws.Range("E5").Value = '  -1.87%  '
$ws.Range("D6").Value = '0.509'
$ws.Range("E6").Value = '  -1.33%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("E8").Value = '  -1.31%  '
$ws.Range("D9").Value = '0.0623'
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("D10").Value = '19.73'
$ws.Range("E10").Value = '  -1.25%  '
$ws.Range("D11").Value = '0.0836'
$ws.Range("E11").Value = '  -1.38%  '
$ws.Range("D12").Value = '1.839.26'
$ws.Range("E12").Value = '  -0.98%  '
$ws.Range("D13").Value = '1.621.25'
$ws.Range("E13").Value = '  -0.45%  '
$ws.Range("D14").Value = '4.08'
$ws.Range("E14").Value = '  -0.66%  '
$ws.Range("E15").Value = '  -1.34%  '
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").Value = '63.97'
$ws.Range("E16").Value = '  -2.63%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '26.809.20'
$ws.Range("E17").Value = '  -0.97%  '
$ws.Range("D18").Value = '0.0₃0732'
$ws.Range("E18").Value = '  -0.04%  '
$ws.Range("D19").Value = '211.21'
$ws.Range("E19").Value = '  -0.98%  '
$ws.Range("E20").Value = '  -0.04%  '
$ws.Range("E21").Value = '  -0.80%  '
$ws.Range("D22").Value = '4.31'
$ws.Range("E22").Value = '  -1.70%  '
$ws.Range("E23").Value = '  -6.80%  '
$ws.Range("E24").Value = '  -1.53%  '
$ws.Range("D25").Value = '146.55'
$ws.Range("E25").Value = '  -0.17%  '
$ws.Range("D26").Value = '7.49'
$ws.Range("E26").Value = '  +1.80%  '
$ws.Range("E27").Value = '  -0.19%  '
$ws.Range("E28").Value = '  -4.14%  '
$ws.Range("E29").Value = '  -1.05%  '
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("E31").Value = '  -1.82%  '
$ws.Range("D32").Value = '3.27'
$ws.Range("E32").Value = '  -2.13%  '
$ws.Range("D33").Value = '0.705'
$ws.Range("E33").Value = '  +30.54%  '
$ws.Range("E34").Value = '  -1.50%  '
$ws.Range("D35").Value = '1.322.93'
$ws.Range("E35").Value = '  +1.80%  '
$ws.Range("D36").Value = '1.54'
$ws.Range("E36").Value = '  -1.21%  '
$ws.Range("D37").Value = '2.44'
$ws.Range("E37").Value = '  -0.32%  '
$ws.Range("E38").Value = '  -1.21%  '
$ws.Range("D39").Value = '0.829'
$ws.Range("E39").Value = '  -1.50%  '
$ws.Range("E40").Value = '  -0.11%  '
$ws.Range("D41").Value = '0.791'
$ws.Range("E41").Value = '  -1.91%  '
$ws.Range("E42").Value = '  -1.57%  '
$ws.Range("D43").Value = '5.29'
$ws.Range("E43").Value = '  -0.22%  '
$ws.Range("D44").Value = '63.53'
$ws.Range("E44").Value = '  +2.59%  '
$ws.Range("D45").Value = '1.751.41'
$ws.Range("E45").Value = '  -1.01%  '
$ws.Range("D46").Value = '89.12'
$ws.Range("E46").Value = '  -1.28%  '
$ws.Range("D47").Value = '1.62'
$ws.Range("E47").Value = '  +2.02%  '
$ws.Range("D48").Value = '0.813'
$ws.Range("E48").Value = '  +4.68%  '
$ws.Range("D49").Value = '0.0512'
$ws.Range("E49").Value = '  +0.05%  '
$ws.Range("D50").Value = '0.0981'
$ws.Range("E50").Value = '  +3.54%  '
$ws.Range("D51").Value = '7.48'
$ws.Range("E51").Value = '  -0.74%  '
